$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# total bazer label + value (added first so it lands first in the shared string table)
$ws.Range("C5").Value = "total bazer"
$ws.Range("C6").Value = 1286

# New header row for the summary table (columns I:L)
$ws.Range("I1").Value = "Supty"
$ws.Range("J1").Value = "Joy"
$ws.Range("K1").Value = "Rana"
$ws.Range("L1").Value = "left"

# New data row for the summary table
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 110
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 214

# Update the selected cell to match target view state
$ws.Range("C8").Select()
